# Refresh the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) with the latest scrape. Every other cell
# (coin name, link, rank index, header, styling) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "97.597.99" },
    @{ Cell = "E2"; Value = "  -1.96%  " },
    @{ Cell = "D3"; Value = "3.424.83" },
    @{ Cell = "E3"; Value = "  +4.03%  " },
    @{ Cell = "E4"; Value = "  +0.06%  " },
    @{ Cell = "D5"; Value = "255.46" },
    @{ Cell = "E5"; Value = "  +0.80%  " },
    @{ Cell = "D6"; Value = "656.08" },
    @{ Cell = "E6"; Value = "  +5.37%  " },
    @{ Cell = "D7"; Value = "1.49" },
    @{ Cell = "E7"; Value = "  +5.01%  " },
    @{ Cell = "D8"; Value = "0.428" },
    @{ Cell = "E8"; Value = "  +7.06%  " },
    @{ Cell = "D9"; Value = "1.06" },
    @{ Cell = "E9"; Value = "  +9.99%  " },
    @{ Cell = "E10"; Value = "  +0.01%  " },
    @{ Cell = "D11"; Value = "3.418.59" },
    @{ Cell = "E11"; Value = "  +3.91%  " },
    @{ Cell = "D12"; Value = "0.213" },
    @{ Cell = "E12"; Value = "  +6.70%  " },
    @{ Cell = "D13"; Value = "42.08" },
    @{ Cell = "E13"; Value = "  +6.87%  " },
    @{ Cell = "D14"; Value = "6.34" },
    @{ Cell = "E14"; Value = "  +15.88%  " },
    @{ Cell = "D15"; Value = "0.0000260" },
    @{ Cell = "E15"; Value = "  +5.28%  " },
    @{ Cell = "D16"; Value = "97.273.25" },
    @{ Cell = "E16"; Value = "  -1.72%  " },
    @{ Cell = "D17"; Value = "4.061.70" },
    @{ Cell = "E17"; Value = "  +4.61%  " },
    @{ Cell = "D18"; Value = "8.70" },
    @{ Cell = "E18"; Value = "  +37.23%  " },
    @{ Cell = "D19"; Value = "3.421.01" },
    @{ Cell = "E19"; Value = "  +4.12%  " },
    @{ Cell = "D20"; Value = "17.60" },
    @{ Cell = "E20"; Value = "  +14.64%  " },
    @{ Cell = "D21"; Value = "0.521" },
    @{ Cell = "E21"; Value = "  +60.66%  " },
    @{ Cell = "D22"; Value = "10.94" },
    @{ Cell = "E22"; Value = "  +17.67%  " },
    @{ Cell = "D23"; Value = "3.47" },
    @{ Cell = "E23"; Value = "  +1.10%  " },
    @{ Cell = "D24"; Value = "508.83" },
    @{ Cell = "E24"; Value = "  +4.26%  " },
    @{ Cell = "D25"; Value = "0.0000206" },
    @{ Cell = "E25"; Value = "  +2.88%  " },
    @{ Cell = "D26"; Value = "6.09" },
    @{ Cell = "E26"; Value = "  +7.92%  " },
    @{ Cell = "D27"; Value = "96.95" },
    @{ Cell = "E27"; Value = "  +8.80%  " },
    @{ Cell = "D28"; Value = "12.80" },
    @{ Cell = "E28"; Value = "  +6.73%  " },
    @{ Cell = "D29"; Value = "3.610.22" },
    @{ Cell = "E29"; Value = "  +5.13%  " },
    @{ Cell = "E30"; Value = "  +14.35%  " },
    @{ Cell = "D31"; Value = "11.46" },
    @{ Cell = "E31"; Value = "  +10.69%  " },
    @{ Cell = "E32"; Value = "  +5.57%  " },
    @{ Cell = "E33"; Value = "  -0.13%  " },
    @{ Cell = "D34"; Value = "0.582" },
    @{ Cell = "E34"; Value = "  +22.91%  " },
    @{ Cell = "E35"; Value = "  +0.11%  " },
    @{ Cell = "D36"; Value = "30.00" },
    @{ Cell = "E36"; Value = "  +7.53%  " },
    @{ Cell = "D37"; Value = "2.28" },
    @{ Cell = "E37"; Value = "  +18.12%  " },
    @{ Cell = "D38"; Value = "7.82" },
    @{ Cell = "E38"; Value = "  +8.52%  " },
    @{ Cell = "D39"; Value = "0.155" },
    @{ Cell = "E39"; Value = "  +4.07%  " },
    @{ Cell = "D40"; Value = "1.43" },
    @{ Cell = "E40"; Value = "  +16.25%  " },
    @{ Cell = "D41"; Value = "518.39" },
    @{ Cell = "E41"; Value = "  +6.24%  " },
    @{ Cell = "D43"; Value = "0.864" },
    @{ Cell = "E43"; Value = "  +12.03%  " },
    @{ Cell = "E44"; Value = "  +1.77%  " },
    @{ Cell = "D45"; Value = "0.0422" },
    @{ Cell = "E45"; Value = "  +26.02%  " },
    @{ Cell = "D46"; Value = "3.30" },
    @{ Cell = "E46"; Value = "  +7.11%  " },
    @{ Cell = "D47"; Value = "5.50" },
    @{ Cell = "E47"; Value = "  +16.62%  " },
    @{ Cell = "D48"; Value = "8.22" },
    @{ Cell = "E48"; Value = "  +13.05%  " },
    @{ Cell = "E49"; Value = "  +0.09%  " },
    @{ Cell = "D50"; Value = "1.58" },
    @{ Cell = "E50"; Value = "  +16.97%  " },
    @{ Cell = "D51"; Value = "2.09" },
    @{ Cell = "E51"; Value = "  +7.30%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $newValue = $u.Value

    # Some "Price" values (column D) are plain decimals (e.g. "1.49"),
    # which Excel's COM layer will happily auto-coerce into a Number on
    # assignment. The source data is text (inlineStr in the OOXML), so we
    # force text entry with a leading apostrophe - exactly like typing
    # '1.49 into the cell - then restore the cell's original Style
    # afterwards, since forcing text entry otherwise stamps a new
    # quote-prefixed style onto the cell.
    $looksNumeric = $newValue -match '^[+-]?\d+(\.\d+)?$'

    if ($looksNumeric) {
        $originalStyle = $cell.Style
        $cell.Value = "'" + $newValue
        $cell.Style = $originalStyle
    } else {
        $cell.Value = $newValue
    }
}
